$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G3").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G5").Value = 'lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G7").Value = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G13").Value = 'ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("G16").Value = 'aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range("G18").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G19").Value = 'Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G21").Value = 'lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G23").Value = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G29").Value = 'ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("G32").Value = 'aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range("G34").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G35").Value = 'eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg'
$ws.Range("G50").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G51").Value = 'eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg'
$ws.Range("G66").Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G67").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G71").Value = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Range("G77").Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G82").Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G83").Value = 'Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G87").Value = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Range("G93").Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G98").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G103").Value = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G114").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G119").Value = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'